# Add a new quarterly sheet "2022-Q3" right after the "总计" summary sheet,
# populate it with the new fund data, and update the summary sheet so the
# "总计" table grows a row (2020-Q4 moves from row 8 to row 9) and gets a
# new row 2 for "2022-Q3".

$wb = $excel.ActiveWorkbook

# Excel constants used below (kept inline so the script has no dependencies
# besides the COM object model):
#   xlCenter           = -4108
#   xlTop               = -4160
#   xlContinuous (line) = 1

# ---------------------------------------------------------------------------
# 1) Insert the new worksheet right after "总计" (position 1) and name it.
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $summarySheet)
$newSheet.Name = "2022-Q3"

function Style-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 2) Fill in the header row (bold, bordered, centered - matching the style
#    used by the header rows of the other quarterly sheets).
# ---------------------------------------------------------------------------
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"
for ($col = 2; $col -le 8; $col++) {
    Style-HeaderCell $newSheet.Cells.Item(1, $col)
}

# ---------------------------------------------------------------------------
# 3) Fill in the two fund data rows. Numeric-looking identifiers/figures are
#    entered with a leading apostrophe so Excel keeps them as text (this
#    preserves things like the leading zero in "014982" and the trailing
#    zero in "0.0860"), then ClearFormats() strips the resulting
#    quote-prefix styling so the cell is plain/unstyled again.
# ---------------------------------------------------------------------------

# Row 2: fund 160416
$newSheet.Cells.Item(2, 1).Value = 0
Style-HeaderCell $newSheet.Cells.Item(2, 1)

$newSheet.Cells.Item(2, 2).Value = "'160416"
$newSheet.Cells.Item(2, 2).ClearFormats()

$newSheet.Cells.Item(2, 3).Value = "华安标普全球石油指数（QDII-LOF）A"

$newSheet.Cells.Item(2, 4).Value = "'2.74"
$newSheet.Cells.Item(2, 4).ClearFormats()

$newSheet.Cells.Item(2, 5).Value = "'93.58"
$newSheet.Cells.Item(2, 5).ClearFormats()

$newSheet.Cells.Item(2, 6).Value = "'3.14"
$newSheet.Cells.Item(2, 6).ClearFormats()

$newSheet.Cells.Item(2, 7).Value = "'0.0860"
$newSheet.Cells.Item(2, 7).ClearFormats()

$newSheet.Cells.Item(2, 8).Value = 7

# Row 3: fund 014982
$newSheet.Cells.Item(3, 1).Value = 1
Style-HeaderCell $newSheet.Cells.Item(3, 1)

$newSheet.Cells.Item(3, 2).Value = "'014982"
$newSheet.Cells.Item(3, 2).ClearFormats()

$newSheet.Cells.Item(3, 3).Value = "华安标普全球石油指数（QDII-LOF）C"

$newSheet.Cells.Item(3, 4).Value = "'0.22"
$newSheet.Cells.Item(3, 4).ClearFormats()

$newSheet.Cells.Item(3, 5).Value = "'93.58"
$newSheet.Cells.Item(3, 5).ClearFormats()

$newSheet.Cells.Item(3, 6).Value = "'3.14"
$newSheet.Cells.Item(3, 6).ClearFormats()

$newSheet.Cells.Item(3, 7).Value = "'0.0069"
$newSheet.Cells.Item(3, 7).ClearFormats()

$newSheet.Cells.Item(3, 8).Value = 7

# ---------------------------------------------------------------------------
# 4) Update the "总计" summary sheet. The table (B/C/D = date, count, value;
#    A = running index) grows from 8 data rows to 9: a new "2022-Q3" row is
#    inserted right under the header and every later row shifts down by one,
#    with "2020-Q4" ending up on row 9. Values are written as literals
#    directly (rather than shifted programmatically) for every row so the
#    saved numbers are exact.
# ---------------------------------------------------------------------------

$summaryRows = @(
    @(2, 0, "2022-Q3", 2, 0.09),
    @(3, 1, "2022-Q2", 3, 0.22),
    @(4, 2, "2022-Q1", 2, 0.13),
    @(5, 3, "2021-Q4", 1, 0.11),
    @(6, 4, "2021-Q3", 2, 0.14),
    @(7, 5, "2021-Q2", 1, 0.14),
    @(8, 6, "2021-Q1", 1, 0.19),
    @(9, 7, "2020-Q4", 1, 0.31)
)

$summarySheet.Cells.Item(2, 1).Value = 0
Style-HeaderCell $summarySheet.Cells.Item(2, 1)
$summarySheet.Cells.Item(2, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 2
$summarySheet.Cells.Item(2, 4).Value = 0.09

$summarySheet.Cells.Item(3, 1).Value = 1
Style-HeaderCell $summarySheet.Cells.Item(3, 1)
$summarySheet.Cells.Item(3, 2).Value = "2022-Q2"
$summarySheet.Cells.Item(3, 3).Value = 3
$summarySheet.Cells.Item(3, 4).Value = 0.22

$summarySheet.Cells.Item(4, 1).Value = 2
Style-HeaderCell $summarySheet.Cells.Item(4, 1)
$summarySheet.Cells.Item(4, 2).Value = "2022-Q1"
$summarySheet.Cells.Item(4, 3).Value = 2
$summarySheet.Cells.Item(4, 4).Value = 0.13

$summarySheet.Cells.Item(5, 1).Value = 3
Style-HeaderCell $summarySheet.Cells.Item(5, 1)
$summarySheet.Cells.Item(5, 2).Value = "2021-Q4"
$summarySheet.Cells.Item(5, 3).Value = 1
$summarySheet.Cells.Item(5, 4).Value = 0.11

$summarySheet.Cells.Item(6, 1).Value = 4
Style-HeaderCell $summarySheet.Cells.Item(6, 1)
$summarySheet.Cells.Item(6, 2).Value = "2021-Q3"
$summarySheet.Cells.Item(6, 3).Value = 2
$summarySheet.Cells.Item(6, 4).Value = 0.14

$summarySheet.Cells.Item(7, 1).Value = 5
Style-HeaderCell $summarySheet.Cells.Item(7, 1)
$summarySheet.Cells.Item(7, 2).Value = "2021-Q2"
$summarySheet.Cells.Item(7, 3).Value = 1
$summarySheet.Cells.Item(7, 4).Value = 0.14

$summarySheet.Cells.Item(8, 1).Value = 6
Style-HeaderCell $summarySheet.Cells.Item(8, 1)
$summarySheet.Cells.Item(8, 2).Value = "2021-Q1"
$summarySheet.Cells.Item(8, 3).Value = 1
$summarySheet.Cells.Item(8, 4).Value = 0.19

$summarySheet.Cells.Item(9, 1).Value = 7
Style-HeaderCell $summarySheet.Cells.Item(9, 1)
$summarySheet.Cells.Item(9, 2).Value = "2020-Q4"
$summarySheet.Cells.Item(9, 3).Value = 1
$summarySheet.Cells.Item(9, 4).Value = 0.31
